$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 23-29 for the MS River stations (USACE source) ---
$newRows = @(
    @{Row=23; Station="01400"; Value=0},
    @{Row=24; Station="01390"; Value=0},
    @{Row=25; Station="01300"; Value=-0.82},
    @{Row=26; Station="01280"; Value=-0.7},
    @{Row=27; Station="01275"; Value=-0.55000000000000004},
    @{Row=28; Station="01260"; Value=-0.76},
    @{Row=29; Station="03780"; Value=-1.1100000000000001}
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = "USACE"

    $ws.Range("B$rowNum").NumberFormat = "@"
    $ws.Range("B$rowNum").Value = $r.Station

    $ws.Range("C$rowNum").NumberFormat = "0.00"
    $ws.Range("C$rowNum").Value = $r.Value
}

# --- Row 12: observation data source changes from USACE to USGS ---
$ws.Range("A12").Value = "USGS"
$ws.Range("B12").Value = "291929089562600"

# --- Selection moves to F14 ---
[void]$ws.Range("F14").Select()
